$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 59

# Columns A-D hold text values (date/time/weekday/week stored as strings,
# matching the existing rows in the sheet). Force text formatting so
# Excel does not auto-convert "2024-01-15" / "14:23:24" / "02" into a
# date, time, or number, then clear the formatting again so no style
# index gets attached to the new cells.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-15"
$ws.Cells.Item($row, 2).Value = "14:23:24"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "02"

$textRange.ClearFormats()

# Columns E-T hold the numeric resale figures for the new entry.
$ws.Cells.Item($row, 5).Value = 138935
$ws.Cells.Item($row, 6).Value = 139029
$ws.Cells.Item($row, 7).Value = 171126
$ws.Cells.Item($row, 8).Value = 148132
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119347
$ws.Cells.Item($row, 11).Value = 220731
$ws.Cells.Item($row, 12).Value = 253947
$ws.Cells.Item($row, 13).Value = 185217
$ws.Cells.Item($row, 14).Value = 110449
$ws.Cells.Item($row, 15).Value = 41123
$ws.Cells.Item($row, 16).Value = 30897
$ws.Cells.Item($row, 17).Value = 73185
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42465
$ws.Cells.Item($row, 20).Value = -1
